$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above current row 20, shifting rows 20-36 down to 21-37
$ws.Rows("20:20").Insert()

# Populate the new row 20 with the new record (reusing style s="2" handled by
# Excel copying formats from the row below on insert, but we set values explicitly)
$ws.Cells.Item(20, 1).Value = 10
$ws.Cells.Item(20, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(20, 3).Value = "La Araucanía"
$ws.Cells.Item(20, 4).Value = 44719
$ws.Cells.Item(20, 5).Value = 9
$ws.Cells.Item(20, 6).Value = 100112010
$ws.Cells.Item(20, 7).Value = "Achicoria"
$ws.Cells.Item(20, 8).Value = "Sin especificar"
$ws.Cells.Item(20, 9).Value = "Primera"
$ws.Cells.Item(20, 10).Value = 50
$ws.Cells.Item(20, 11).Value = 10000
$ws.Cells.Item(20, 12).Value = 10000
$ws.Cells.Item(20, 13).Value = 10000
$ws.Cells.Item(20, 14).Value = "`$/caja 18 unidades"
$ws.Cells.Item(20, 15).Value = "Región Metropolitana"
$ws.Cells.Item(20, 16).Value = 556
$ws.Cells.Item(20, 17).Value = 18
$ws.Cells.Item(20, 18).Value = "Hortaliza"
